$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 329.83334
$ws.Range("I12").Value = 361.8
$ws.Range("K12").Value = 361.8
$ws.Range("M12").Value = -191.8
$ws.Range("H51").Value = 7184.364
$ws.Range("I51").Value = 6756.6
$ws.Range("K51").Value = 6756.6
$ws.Range("M51").Value = -6272.6
$ws.Range("H53").Value = 411.44446
$ws.Range("J53").Value = 444.125
$ws.Range("L53").Value = 444.125
$ws.Range("N53").Value = -1718.125
$ws.Range("I113").Value = 11113564
$ws.Range("K113").Value = 11113564
$ws.Range("M113").Value = -11110310
$ws.Range("H120").Value = 53190
$ws.Range("J120").Value = 53190
$ws.Range("L120").Value = 53190
$ws.Range("N120").Value = -62866
$ws.Range("H135").Value = 1004.7
$ws.Range("I135").Value = 1005.0357
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9045.3213
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6510.3213
$ws.Range("N135").Value = -14070
$ws.Range("H138").Value = 2362.76
$ws.Range("I138").Value = 1221
$ws.Range("J138").Value = 2648.2
$ws.Range("K138").Value = 3663
$ws.Range("L138").Value = 7944.599999999999
$ws.Range("M138").Value = 1477
$ws.Range("N138").Value = -18224.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 29011.334
$ws.Range("J37").Value = 45000
$ws.Range("L37").Value = 45000
$ws.Range("N37").Value = -45546
$ws.Range("H45").Value = 1952
$ws.Range("I45").Value = 1336.0714
$ws.Range("K45").Value = 1336.0714
$ws.Range("M45").Value = -959.0714
$ws.Range("H97").Value = 1115.08
$ws.Range("I97").Value = 1115.08
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1115.08
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -619.0799999999999
$ws.Range("N97").ClearContents()
$ws.Range("H112").Value = 36981.332
$ws.Range("J112").Value = 36981.332
$ws.Range("L112").Value = 36981.332
$ws.Range("N112").Value = -39935.332
$ws.Range("H122").Value = 1261.1904
$ws.Range("I122").Value = 1082.6111
$ws.Range("K122").Value = 3247.8333
$ws.Range("M122").Value = -797.8333000000002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29516
$ws.Range("H134").Value = 67122.625
$ws.Range("I134").Value = 844.4545000000001
$ws.Range("J134").Value = 212934.6
$ws.Range("K134").Value = 2533.3635
$ws.Range("L134").Value = 638803.8
$ws.Range("M134").Value = 1.636499999999614
$ws.Range("N134").Value = -643873.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6027.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 6027.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14542.846
$ws.Range("I3").Value = 10999.2
$ws.Range("K3").Value = 32997.60000000001
$ws.Range("M3").Value = -32885.60000000001
$ws.Range("H99").Value = 2841.6667
$ws.Range("I99").Value = 2841.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8525.000100000001
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6279.000100000001
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 502.91306
$ws.Range("I107").Value = 382.27274
$ws.Range("J107").Value = 613.5
$ws.Range("K107").Value = 1146.81822
$ws.Range("L107").Value = 1840.5
$ws.Range("M107").Value = 773.1817799999999
$ws.Range("N107").Value = -5680.5
$ws.Range("H122").Value = 741.1667
$ws.Range("J122").Value = 2091.875
$ws.Range("L122").Value = 18826.875
$ws.Range("N122").Value = -23726.875
$ws.Range("H132").Value = 1246.9131
$ws.Range("I132").Value = 1216
$ws.Range("J132").Value = 1275.25
$ws.Range("K132").Value = 10944
$ws.Range("L132").Value = 11477.25
$ws.Range("M132").Value = -8414
$ws.Range("N132").Value = -16537.25
$ws.Range("H133").Value = 6018.1816
$ws.Range("I133").Value = 5200
$ws.Range("J133").Value = 7000
$ws.Range("K133").Value = 15600
$ws.Range("L133").Value = 21000
$ws.Range("M133").Value = -10540
$ws.Range("N133").Value = -31120
$ws.Range("H140").Value = 160196.42
$ws.Range("I140").Value = 168912.89
$ws.Range("K140").Value = 506738.67
$ws.Range("M140").Value = -501558.67
$ws.Range("H141").Value = 223627.72
$ws.Range("I141").Value = 603358
$ws.Range("J141").Value = 12666.444
$ws.Range("K141").Value = 1810074
$ws.Range("L141").Value = 37999.33199999999
$ws.Range("M141").Value = -1804894
$ws.Range("N141").Value = -48359.33199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2375.842
$ws.Range("I107").Value = 1689
$ws.Range("J107").Value = 3553.2856
$ws.Range("K107").Value = 1689
$ws.Range("L107").Value = 3553.2856
$ws.Range("M107").Value = 231
$ws.Range("N107").Value = -7393.2856
$ws.Range("H111").Value = 58935.75
$ws.Range("J111").Value = 58935.75
$ws.Range("L111").Value = 58935.75
$ws.Range("N111").Value = -65069.75
$ws.Range("H122").Value = 2332.5
$ws.Range("I122").Value = 1998.75
$ws.Range("K122").Value = 5996.25
$ws.Range("M122").Value = -3546.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18277864
$ws.Range("I7").Value = 50003624
$ws.Range("J7").Value = 148857.72
$ws.Range("K7").Value = 50003624
$ws.Range("L7").Value = 148857.72
$ws.Range("M7").Value = -50003512
$ws.Range("N7").Value = -149081.72
$ws.Range("H74").Value = 67550
$ws.Range("J74").Value = 67550
$ws.Range("L74").Value = 67550
$ws.Range("N74").Value = -69546
$ws.Range("H77").Value = 67550
$ws.Range("J77").Value = 67550
$ws.Range("L77").Value = 202650
$ws.Range("N77").Value = -212634
$ws.Range("H93").Value = 76926010
$ws.Range("I93").Value = 90911784
$ws.Range("K93").Value = 90911784
$ws.Range("M93").Value = -90910536
$ws.Range("H100").Value = 2421.1
$ws.Range("I100").Value = 2201.5
$ws.Range("J100").Value = 3299.5
$ws.Range("K100").Value = 2201.5
$ws.Range("L100").Value = 3299.5
$ws.Range("M100").Value = -1660.5
$ws.Range("N100").Value = -4381.5
$ws.Range("H122").Value = 4666.6772
$ws.Range("I122").Value = 4382.15
$ws.Range("K122").Value = 13146.45
$ws.Range("M122").Value = -10696.45
$ws.Range("H126").Value = 18277864
$ws.Range("I126").Value = 50003624
$ws.Range("J126").Value = 148857.72
$ws.Range("K126").Value = 150010872
$ws.Range("L126").Value = 446573.16
$ws.Range("M126").Value = -150008402
$ws.Range("N126").Value = -451513.16
$ws.Range("H127").Value = 121442.25
$ws.Range("J127").Value = 121442.25
$ws.Range("L127").Value = 121442.25
$ws.Range("N127").Value = -131362.25
$ws.Range("H132").Value = 345315.75
$ws.Range("I132").Value = 324196.75
$ws.Range("K132").Value = 972590.25
$ws.Range("M132").Value = -970060.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 45987.145
$ws.Range("I52").Value = 48283.4
$ws.Range("J52").Value = 40246.5
$ws.Range("K52").Value = 48283.4
$ws.Range("L52").Value = 40246.5
$ws.Range("M52").Value = -48057.4
$ws.Range("N52").Value = -40698.5
$ws.Range("H62").Value = 5889275
$ws.Range("J62").Value = 25005412
$ws.Range("L62").Value = 25005412
$ws.Range("N62").Value = -25006660
$ws.Range("H65").Value = 5889275
$ws.Range("J65").Value = 25005412
$ws.Range("L65").Value = 125027060
$ws.Range("N65").Value = -125033300
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 5612.091
$ws.Range("I122").Value = 3331.16
$ws.Range("K122").Value = 9993.48
$ws.Range("M122").Value = -7543.48
$ws.Range("H126").Value = 2531.111
$ws.Range("I126").Value = 2279.1667
$ws.Range("K126").Value = 6837.500100000001
$ws.Range("M126").Value = -4367.500100000001
$ws.Range("H132").Value = 1559.675
$ws.Range("I132").Value = 1357.8611
$ws.Range("K132").Value = 4073.5833
$ws.Range("M132").Value = -1543.5833
